# Remove the post row for "「こわれている」مكسر" (row 735) and let Excel
# shift all subsequent rows up by one, as in the source commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(735).Delete()
